$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Nama" column (column B) is no longer needed in the import template,
# so select it (as a user would before deleting) and remove it entirely.
# This shifts Nomor_SP2DK, Tanggal_SP2DK, Tahun and Potensi Awal one
# column to the left (C->B, D->C, E->D, F->E).
$ws.Columns.Item(2).Select()
$ws.Columns.Item(2).Delete()
